# Updated symbol list on Fri Dec 30 19:49:25 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column values are written back as plain text (not
# auto-converted to numbers) so formatting such as trailing zeros and
# exact decimal representation is preserved.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D19","D20","D22","D23","D26","D27","D40","D41","D42","D43","D44","D45","D47","D48")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Row 2 - BNB
$ws.Range("D2").Value = "245.01"

# Row 3 - OKB
$ws.Range("D3").Value = "25.17"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "4.980"

# Row 5 - Cronos
$ws.Range("D5").Value = "0.05616"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "6.546"

# Row 7 - GateToken
$ws.Range("D7").Value = "3.006"

# Row 8 - MXToken
$ws.Range("D8").Value = "0.8116"

# Row 9 - FTXToken
$ws.Range("D9").Value = "0.8398"

# Row 10 - WazirX
$ws.Range("D10").Value = "0.1336"

# Row 11 - was LiechtensteinCryptoassetsExchange, now MandalaExchangeToken
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.06938"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12 - BitrueCoin
$ws.Range("D12").Value = "0.02846"

# Row 13 - BitMartToken
$ws.Range("D13").Value = "0.09409"

# Row 14 - BitForexToken
$ws.Range("D14").Value = "0.001515"

# Row 15 - One
$ws.Range("D15").Value = "0.0005979"
$ws.Range("E15").Value = "14OneONEWorstin24h"

# Row 16 - TigerCash
$ws.Range("D16").Value = "0.006252"

# Row 19 - BitpandaEcosystemToken
$ws.Range("D19").Value = "0.3196"

# Row 20 - was MandalaExchangeToken, now LiechtensteinCryptoassetsExchange
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "0.03250"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"

# Row 22
$ws.Range("D22").Value = "3.746"

# Row 23
$ws.Range("D23").Value = "0.04678"

# Row 26
$ws.Range("D26").Value = "0.004528"

# Row 27
$ws.Range("D27").Value = "0.00009697"

# Row 40
$ws.Range("D40").Value = "0.03643"

# Row 41 - BKEXToken
$ws.Range("D41").Value = "0.1366"

# Row 42 - was CEJI, now KickToken
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "0.006248"
$ws.Range("E42").Value = "41KickTokenKICK"

# Row 43 - was KickToken, now CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002723"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "0.008061"

# Row 45
$ws.Range("D45").Value = "0.00005272"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "0.1800"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Row 48 - BOLO
$ws.Range("D48").Value = "0.002043"
